$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -6
$ws.Range("F3").Value = -3
$ws.Range("F4").Value = -4
$ws.Range("F5").Value = -2
$ws.Range("F6").Value = 2
$ws.Range("F8").Value = -3
$ws.Range("F12").Value = 4
$ws.Range("F15").Value = 1
$ws.Range("F19").Value = -7
$ws.Range("F20").Value = -2
$ws.Range("F21").Value = -3
$ws.Range("F26").Value = 1
$ws.Range("F33").Value = 4
$ws.Range("F38").Value = 4
$ws.Range("F43").Value = -1
$ws.Range("F44").Value = -5
$ws.Range("F47").Value = -5
$ws.Range("F48").Value = -2
$ws.Range("F53").Value = 3
$ws.Range("F61").Value = 1
$ws.Range("F62").Value = -6
$ws.Range("F66").Value = 1
$ws.Range("F70").Value = 0
$ws.Range("F71").Value = 4
$ws.Range("F72").Value = -3
